# Auto-generated script applying 211 cell value updates across 8 worksheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the authoritative diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 190
$ws.Range("I9").Value = 128.4
$ws.Range("J9").Value = 258.44446
$ws.Range("K9").Value = 128.4
$ws.Range("L9").Value = 258.44446
$ws.Range("M9").Value = 40.59999999999999
$ws.Range("N9").Value = -596.4444599999999
$ws.Range("H33").Value = 238.44444
$ws.Range("J33").Value = 176
$ws.Range("L33").Value = 176
$ws.Range("N33").Value = -634
$ws.Range("H74").Value = 6226.1333
$ws.Range("I74").Value = 6226.1333
$ws.Range("K74").Value = 6226.1333
$ws.Range("M74").Value = -5290.1333
$ws.Range("H77").Value = 6226.1333
$ws.Range("I77").Value = 6226.1333
$ws.Range("K77").Value = 31130.6665
$ws.Range("M77").Value = -26450.6665
$ws.Range("H138").Value = 4447.636
$ws.Range("I138").Value = 2844.818
$ws.Range("J138").Value = 4981.909
$ws.Range("K138").Value = 8534.454000000002
$ws.Range("L138").Value = 14945.727
$ws.Range("M138").Value = -3394.454000000002
$ws.Range("N138").Value = -25225.727

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1763.8823
$ws.Range("I2").Value = 1691.4333
$ws.Range("K2").Value = 1691.4333
$ws.Range("M2").Value = -1578.4333
$ws.Range("H61").Value = 1291.3
$ws.Range("I61").Value = 1114.125
$ws.Range("J61").Value = 2000
$ws.Range("K61").Value = 1114.125
$ws.Range("L61").Value = 2000
$ws.Range("M61").Value = -902.125
$ws.Range("N61").Value = -2424
$ws.Range("H74").Value = 1330.35
$ws.Range("I74").Value = 1330.35
$ws.Range("K74").Value = 1330.35
$ws.Range("M74").Value = -456.3499999999999
$ws.Range("H77").Value = 1330.35
$ws.Range("I77").Value = 1330.35
$ws.Range("K77").Value = 6651.75
$ws.Range("M77").Value = -2283.75
$ws.Range("H88").Value = 3487.9092
$ws.Range("I88").Value = 2334.3333
$ws.Range("J88").Value = 3920.5
$ws.Range("K88").Value = 2334.3333
$ws.Range("L88").Value = 3920.5
$ws.Range("M88").Value = -1928.3333
$ws.Range("N88").Value = -4732.5
$ws.Range("H91").Value = 3487.9092
$ws.Range("I91").Value = 2334.3333
$ws.Range("J91").Value = 3920.5
$ws.Range("K91").Value = 2334.3333
$ws.Range("L91").Value = 3920.5
$ws.Range("M91").Value = -930.3332999999998
$ws.Range("N91").Value = -6728.5
$ws.Range("H102").Value = 4950.3
$ws.Range("I102").Value = 4176.4287
$ws.Range("K102").Value = 4176.4287
$ws.Range("M102").Value = -2554.4287
$ws.Range("H116").Value = 1763.8823
$ws.Range("I116").Value = 1691.4333
$ws.Range("K116").Value = 1691.4333
$ws.Range("M116").Value = 602.5667000000001
$ws.Range("H136").Value = 1291.3
$ws.Range("I136").Value = 1114.125
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 3342.375
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -792.375
$ws.Range("N136").Value = -11100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1763.8823
$ws.Range("I3").Value = 1691.4333
$ws.Range("K3").Value = 1691.4333
$ws.Range("M3").Value = -1577.4333
$ws.Range("H94").Value = 3042.3333
$ws.Range("J94").Value = 4775.857
$ws.Range("L94").Value = 4775.857
$ws.Range("N94").Value = -5677.857
$ws.Range("H99").Value = 4990.3335
$ws.Range("I99").Value = 4856.5386
$ws.Range("J99").Value = 5338.2
$ws.Range("K99").Value = 4856.5386
$ws.Range("L99").Value = 5338.2
$ws.Range("M99").Value = -3358.5386
$ws.Range("N99").Value = -8334.200000000001
$ws.Range("H134").Value = 4508.5
$ws.Range("I134").Value = 3433.4546
$ws.Range("K134").Value = 10300.3638
$ws.Range("M134").Value = -7765.363799999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 825.6667
$ws.Range("I3").Value = 825.6667
$ws.Range("K3").Value = 825.6667
$ws.Range("M3").Value = -712.6667
$ws.Range("H22").Value = 1527
$ws.Range("I22").Value = 986
$ws.Range("J22").Value = 3150
$ws.Range("K22").Value = 986
$ws.Range("L22").Value = 3150
$ws.Range("M22").Value = -636
$ws.Range("N22").Value = -3850
$ws.Range("H62").Value = 90913120
$ws.Range("I62").Value = 250003740
$ws.Range("J62").Value = 4195.7144
$ws.Range("K62").Value = 250003740
$ws.Range("L62").Value = 4195.7144
$ws.Range("M62").Value = -250003116
$ws.Range("N62").Value = -5443.7144
$ws.Range("H65").Value = 90913120
$ws.Range("I65").Value = 250003740
$ws.Range("J65").Value = 4195.7144
$ws.Range("K65").Value = 1250018700
$ws.Range("L65").Value = 20978.572
$ws.Range("M65").Value = -1250015580
$ws.Range("N65").Value = -27218.572
$ws.Range("H94").Value = 3224.7058
$ws.Range("J94").Value = 3311.5
$ws.Range("L94").Value = 3311.5
$ws.Range("N94").Value = -4213.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 45.53846
$ws.Range("I12").Value = 46.272728
$ws.Range("J12").Value = 45
$ws.Range("K12").Value = 138.818184
$ws.Range("L12").Value = 135
$ws.Range("M12").Value = 34.181816
$ws.Range("N12").Value = -481
$ws.Range("H14").Value = 11229
$ws.Range("I14").Value = 11229
$ws.Range("K14").Value = 33687
$ws.Range("M14").Value = -33514
$ws.Range("H68").Value = 2059.8667
$ws.Range("J68").Value = 2507.182
$ws.Range("L68").Value = 7521.545999999999
$ws.Range("N68").Value = -9143.545999999998
$ws.Range("H71").Value = 2059.8667
$ws.Range("J71").Value = 2507.182
$ws.Range("L71").Value = 22564.638
$ws.Range("N71").Value = -30676.638
$ws.Range("H107").Value = 604.45
$ws.Range("I107").Value = 463.66666
$ws.Range("J107").Value = 719.63635
$ws.Range("K107").Value = 1390.99998
$ws.Range("L107").Value = 2158.90905
$ws.Range("M107").Value = 529.0000199999999
$ws.Range("N107").Value = -5998.90905
$ws.Range("H119").Value = 4267.5
$ws.Range("I119").Value = 1134.6
$ws.Range("K119").Value = 3403.8
$ws.Range("M119").Value = 1434.2
$ws.Range("H134").Value = 3000
$ws.Range("I134").Value = 3000
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 9000
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -3930
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 6711.7407
$ws.Range("I102").Value = 7658.85
$ws.Range("K102").Value = 7658.85
$ws.Range("M102").Value = -6036.85

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2131.1667
$ws.Range("I22").Value = 2358.4
$ws.Range("J22").Value = 995
$ws.Range("K22").Value = 2358.4
$ws.Range("L22").Value = 995
$ws.Range("M22").Value = -2063.4
$ws.Range("N22").Value = -1585
$ws.Range("H27").Value = 2131.1667
$ws.Range("I27").Value = 2358.4
$ws.Range("J27").Value = 995
$ws.Range("K27").Value = 2358.4
$ws.Range("L27").Value = 995
$ws.Range("M27").Value = -2251.4
$ws.Range("N27").Value = -1209
$ws.Range("H93").Value = 3849.5737
$ws.Range("I93").Value = 4064.6511
$ws.Range("J93").Value = 3335.7778
$ws.Range("K93").Value = 4064.6511
$ws.Range("L93").Value = 3335.7778
$ws.Range("M93").Value = -2816.6511
$ws.Range("N93").Value = -5831.7778

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 11108
$ws.Range("I12").Value = 12897
$ws.Range("J12").Value = 8424.5
$ws.Range("K12").Value = 12897
$ws.Range("L12").Value = 8424.5
$ws.Range("M12").Value = -12755
$ws.Range("N12").Value = -8708.5
$ws.Range("H62").Value = 4777.5
$ws.Range("I62").Value = 4777.5
$ws.Range("K62").Value = 4777.5
$ws.Range("M62").Value = -4153.5
$ws.Range("H65").Value = 4777.5
$ws.Range("I65").Value = 4777.5
$ws.Range("K65").Value = 23887.5
$ws.Range("M65").Value = -20767.5
$ws.Range("H113").Value = 1712.875
$ws.Range("I113").Value = 1533.625
$ws.Range("K113").Value = 4600.875
$ws.Range("M113").Value = -2430.875
$ws.Range("H124").Value = 103761
$ws.Range("J124").Value = 103761
$ws.Range("L124").Value = 103761
$ws.Range("N124").Value = -113581
$ws.Range("H132").Value = 8216.950999999999
$ws.Range("I132").Value = 6479.5864
$ws.Range("K132").Value = 19438.7592
$ws.Range("M132").Value = -16908.7592

Write-Output "Applied 211 cell updates"